$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.155900359153748
$ws.Range("B1").Value = 2.396236658096313
$ws.Range("D1").Value = 2.384890556335449
$ws.Range("E1").Value = 1.225401878356934
